$wb = $excel.ActiveWorkbook

# --- 1. Update the "bulk storage" explanation text to mention the new [ShowBlip] flag ---
$wsField = $wb.Worksheets.Item("Field Explanation")
$wsField.Range("C9").Value = "Used as 'bulk storage' for the boolean variables [Ownable][Owned][ContextMission][ShowBlip]"

# --- 2. Update the Flags sample values on the "Properties Table" sheet to include the new flag digit ---
$wsProps = $wb.Worksheets.Item("Properties Table")
$wsProps.Range("C2").Value = 1011
$wsProps.Range("C3").Value = 1002
$wsProps.Range("C4").Value = 1001
$wsProps.Range("C5").Value = 1011
$wsProps.Range("C6").Value = 1001
$wsProps.Range("C7").Value = 1001
$wsProps.Range("C8").Value = 1001
$wsProps.Range("C9").Value = 1001
$wsProps.Range("C10").Value = 1001
$wsProps.Range("C11").Value = 1001

# --- 3. Move the active tab / selection from "Field Explanation" to "Properties Table" ---
$wsField.Range("C9").Select() | Out-Null
$wsProps.Activate() | Out-Null
$wsProps.Range("B13").Select() | Out-Null
